$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grow the table from 2 columns (Services, Amount) to 5 columns ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E2"))

# --- Set the new header row (table columns + shared strings follow automatically) ---
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "ShortCode"
$ws.Range("C1").Value = "Group"
$ws.Range("D1").Value = "SubGroup"
$ws.Range("E1").Value = "Amount"

# --- New header cells (C1:E1) need to carry the same header style as A1/B1 ---
$ws.Range("A1").Copy()
$ws.Range("C1:E1").PasteSpecial(-4122)

# --- Column widths to roughly match the real "best fit" sizing ---
# (the headless engine only stores widths on a 1/6-character grid, so these
#  inputs are chosen to land on the closest achievable grid point to the
#  true Excel auto-fit pixel widths of 79px/93px/84px/75px)
$ws.Columns.Item(1).ColumnWidth = 10.5
$ws.Columns.Item(2).ColumnWidth = 12.5
$ws.Columns.Item(4).ColumnWidth = 11.15
$ws.Columns.Item(5).ColumnWidth = 9.83

# --- Move the active selection ---
$ws.Range("G4").Select() | Out-Null
